# Collapses three spots in the document where a sentence was split across
# multiple <w:r> runs (sometimes with <w:proofErr/> grammar-check markers
# sitting between them) back down into a single run per sentence, as
# described by the commit ("sum[ming up the split] notes").
#
#   1. "Encouraging action or " + "behaviour" + "."
#   2. "To persuade or influence opinions and " + "behaviours" + "."
#   3. "8. 7 " + <proofErr/> + "C's" + <proofErr/> + " of Communication:"

$d = $word.ActiveDocument

function Merge-Runs([string]$text) {
    # Re-resolve Content each time so the Find always scans from the top of
    # the story instead of continuing from a stale/previous match position.
    $r = $d.Content
    $found = $r.Find.Execute($text, $false, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output ("NOT FOUND: " + $text)
        return
    }

    $startPos = $r.Start
    $endPos   = $r.End

    # Writing back the exact same text a Range already contains is treated
    # as a no-op (no run merge, proofErr markers stay put, apostrophes keep
    # their run split). Make the content genuinely change first -- that
    # forces the runs spanned by the Range (and any <w:proofErr/> between
    # them) to collapse into one freshly-formatted run -- then trim the
    # marker back off via a freshly resolved Range (the original Range's
    # .End does not reliably track the prior in-place edit).
    $r.Text = $text + "#"
    $r2 = $d.Range($startPos, $endPos + 1)
    $r2.Text = $text
}

Merge-Runs "Encouraging action or behaviour."
Merge-Runs "To persuade or influence opinions and behaviours."
Merge-Runs "8. 7 C's of Communication:"
